# Read & Write Vector Capability
# - Adds defined names covering whole ranges (vectors) in addition to the
#   existing per-component defined names, for HOUSING moments and products
#   of inertia.
# - Updates the AIR_NUT.point sample vector (C9:E9) with new test values.
# - Rewrites the HOUSING MOI values (C17:C19) with new (smaller) sample data.
# - Adds a new "HOUSING POI" (products of inertia) block (rows 21-24) with
#   labels PXY/PXZ/PYZ and their values.
# - Updates the view so the new POI block is visible/selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gearbox Tests")

# --- New "vector" defined names, referring to the whole component ranges ---
$wb.Names.Add("HOUSING.moments_of_inertia_centroidal", "='Gearbox Tests'!`$C`$17:`$C`$19")
$wb.Names.Add("HOUSING.products_of_inertia_centroidal", "='Gearbox Tests'!`$C`$22:`$C`$24")

# --- AIR_NUT.point sample vector update (row 9) ---
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 2.11
$ws.Range("E9").Value = 9.99

# --- HOUSING MOI values update (rows 17-19) ---
$ws.Range("C17").Value = 0.012
$ws.Range("C18").Value = 0.11
$ws.Range("C19").Value = 0.99

# Keep these rows at the standard row height explicitly.
$ws.Rows.Item(17).RowHeight = 14.5
$ws.Rows.Item(18).RowHeight = 14.5
$ws.Rows.Item(19).RowHeight = 14.5

# --- New HOUSING POI block (rows 21-24) ---
$ws.Range("B21").Value = "HOUSING POI"

$ws.Range("B22").Value = "PXY"
$ws.Range("C22").Value = -8350.6714059345395

$ws.Range("B23").Value = "PXZ"
$ws.Range("C23").Value = -48498.809284873329

$ws.Range("B24").Value = "PYZ"
$ws.Range("C24").Value = 26209.492442448263

# --- View: scroll/zoom to the new block and select it ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 130
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C22:C24").Select()
